# Daily attendance processing - 2025-12-16 14:32:10
# Swap the order of the first two comma-separated names in the "Recorded By"
# (column G) cell of the "Session Analysis Results" sheet, wherever that
# cell's value is one of the known "System"/"admin@admin.com"/"dnasr281@gmail.com"
# combinations recorded by the attendance-processing job. Single-author cells
# (e.g. just "System" or just "dnasr281@gmail.com") and the
# "System, backup@backdoor.com" pairing are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Text

    if ($current -eq $null -or $current -eq "") {
        continue
    }

    $parts = $current -split ", "
    if ($parts.Length -lt 2) {
        continue
    }

    if ($current -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($current -eq "admin@admin.com, System") {
        $cell.Value = "System, admin@admin.com"
    }
    elseif ($current -eq "admin@admin.com, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, admin@admin.com"
    }
    elseif ($current -eq "system, System, backup@backdoor.com") {
        $cell.Value = "System, system, backup@backdoor.com"
    }
}
